$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$B2 = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE g.library_strategy in ['Archer Fusion']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id limit 100
'@
$B3 = @'
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.library_strategy in ['Archer Fusion']
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(samp.sample_tumor_status,'') as `Tumor`,
    coalesce(samp.sample_type,'') as `Analyte Type`
ORDER BY samp.sample_id limit 100
'@
$B4 = @'
Match (f)<--(g:genomic_info)
WHERE g.library_strategy in ['Archer Fusion']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag
WHERE g.library_strategy in ['Archer Fusion']
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name,'') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id, '') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER BY f.file_name limit 100
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@
$C = @'
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.library_strategy in ['Archer Fusion']
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.library_strategy in ['Archer Fusion']
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag
    WHERE g.library_strategy in ['Archer Fusion']
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS `Files`
'@

# Row 2: ParticipantsTab
$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("B2").Value = $B2
$ws.Range("C2").Value = $C
$ws.Range("D2").Value = "TC03_CDS_Filter_LibraryStrategy-ArcherFusion_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC03_CDS_Filter_LibraryStrategy-ArcherFusion_WebData.xlsx"

# Row 3: SamplesTab
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $B3
$ws.Range("C3").Value = $C
$ws.Range("D3").Value = "TC03_CDS_Filter_LibraryStrategy-ArcherFusion_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC03_CDS_Filter_LibraryStrategy-ArcherFusion_WebData.xlsx"

# Row 4: FilesTab
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $B4
$ws.Range("C4").Value = $C
$ws.Range("D4").Value = "TC03_CDS_Filter_LibraryStrategy-ArcherFusion_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC03_CDS_Filter_LibraryStrategy-ArcherFusion_WebData.xlsx"

# Preserve the original row heights (explicit custom heights, unchanged by this edit)
$ws.Rows.Item(2).RowHeight = 242.25
$ws.Rows.Item(3).RowHeight = 260.25
$ws.Rows.Item(4).RowHeight = 279.75

# Update the sheet view selection (was B2, now D3); also drop the topLeftCell freeze-to B1
$ws.Range("D3").Select()
